# Replace the header row (B1:P1) and the "trial #" column (A2:A10),
# which previously held shared-string text labels ("x1".."x9","y1".."y15"),
# with a numeric sequence generated by formulas:
#   B1 = 0 (literal)
#   C1 = =1+B1                     (typed formula)
#   D1:P1 = =1+C1                  (fill right -> shared formula group)
#   A2 = 0 (literal)
#   A3 = =1+A2                     (typed formula)
#   A4:A10 = =1+A3                 (fill down -> shared formula group)
# This removes every reference to the shared strings table, leaving it
# empty, and leaves the active selection on A3:A10 (the range that was
# just filled down).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: B1 literal 0, then C1.. filled right with "=1+<prev>" ---
$ws.Range("B1").Value = 0
$ws.Range("C1").Formula = "=1+B1"
$ws.Range("D1:P1").Formula = "=1+C1"

# --- Column A: A2 literal 0, then A3.. filled down with "=1+<prev>" ---
$ws.Range("A2").Value = 0
$ws.Range("A3").Formula = "=1+A2"
$ws.Range("A4:A10").Formula = "=1+A3"

# Leave the selection where it would be after the fill-down of A3:A10.
$ws.Range("A3:A10").Select()
